{"js": "// Office.js (Word JavaScript API) script\n// 1) Rewrite the first body paragraph so \"...ecological systems,\" becomes\n//    \"...ecological systems (Poisot, Stouffer, and K\u00e9fi 2016),\" i.e. insert\n//    the in-text citation right after \"ecological systems\" (the comma that\n//    used to directly follow \"systems\" now trails the citation instead).\n// 2) Append a new \"Bibliography\"-styled paragraph after the References\n//    heading with the full reference entry (including an italic journal\n//    title and a hyperlinked DOI), wrapped in \"refs\"/\"ref-PoisStou16\"\n//    bookmarks.\n\nconst body = context.document.body;\n\n// --- 1. Update the introduction sentence ------------------------------\nconst hits = body.search(\"ecological systems,\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  hits.items[0].insertText(\n    \"ecological systems (Poisot, Stouffer, and K\u00e9fi 2016),\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// --- 2. Append the bibliography entry ----------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst referencesHeading = paragraphs.items[paragraphs.items.length - 1];\nconst newPara = referencesHeading.insertParagraph(\"\", \"After\");\nnewPara.style = \"Bibliography\";\nawait context.sync();\n\n// Leading citation text (authors, year, title).\nnewPara.insertText(\n  \"Poisot, Timoth\u00e9e, Daniel B. Stouffer, and Sonia K\u00e9fi. 2016. \u201cDescribe, \" +\n    \"Understand and Predict: Why Do We Need Networks in Ecology?\u201d\",\n  \"End\"\n);\nawait context.sync();\n\nnewPara.insertText(\" \", \"End\");\nawait context.sync();\n\n// Italic journal title.\nconst journalRange = newPara.insertText(\"Functional Ecology\", \"End\");\njournalRange.font.italic = true;\nawait context.sync();\n\nnewPara.insertText(\" \", \"End\");\nawait context.sync();\n\n// Volume/issue/page range.\nnewPara.insertText(\"30 (12): 1878\u201382.\", \"End\");\nawait context.sync();\n\nnewPara.insertText(\" \", \"End\");\nawait context.sync();\n\n// Hyperlinked DOI.\nconst doiUrl = \"https://doi.org/10.1111/1365-2435.12799\";\nconst linkRange = newPara.insertText(doiUrl, \"End\");\nlinkRange.hyperlink = doiUrl;\nawait context.sync();\n\n// Trailing period.\nnewPara.insertText(\".\", \"End\");\nawait context.sync();\n\n// --- 3. Wrap the new paragraph in bookmarks -----------------------------\n// Inserted innermost-first so the stack-like insertBookmark ordering\n// produces \"refs\" as the outer bookmark and \"ref-PoisStou16\" as the inner\n// one, matching the target start/end nesting order.\nconst bookmarkRange1 = newPara.getRange();\nbookmarkRange1.insertBookmark(\"ref-PoisStou16\");\nawait context.sync();\n\nconst bookmarkRange2 = newPara.getRange();\nbookmarkRange2.insertBookmark(\"refs\");\nawait context.sync();\n", "ps1": "# Word COM interop script\n# 1) Rewrite the first body paragraph so the sentence \"...ecological\n#    systems,\" becomes \"...ecological systems (Poisot, Stouffer, and\n#    Kefi 2016),\" i.e. insert the in-text citation right after\n#    \"ecological systems\" (dropping the comma that used to follow it,\n#    since the comma now trails the citation instead).\n# 2) Append a new \"Bibliography\"-styled paragraph after the References\n#    heading with the full reference entry (including an italic journal\n#    title and a hyperlinked DOI), wrapped in \"refs\"/\"ref-PoisStou16\"\n#    bookmarks.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the introduction sentence -----------------------------\n$find = $d.Content\n$found = $find.Find.Execute(\"ecological systems,\")\nif ($found) {\n    $find.Text = \"ecological systems (Poisot, Stouffer, and K\u00e9fi 2016),\"\n}\n\n# --- 2. Append the bibliography entry --------------------------------\n$lastPara = $d.Paragraphs.Last\n$endOfDoc = $lastPara.Range\n$endOfDoc.Collapse(0)               # wdCollapseEnd\n$endOfDoc.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Last\n$newPara.Style = \"Bibliography\"\n\n# Add the hyperlink first (into the still-empty paragraph) -- this\n# engine's Hyperlinks.Add always lands its text at the start of the\n# paragraph range supplied, so everything else gets threaded in with\n# InsertBefore/InsertAfter relative to a freshly-fetched paragraph\n# range instead of relying on the hyperlink's own position.\n$linkRng = $newPara.Range\n$linkRng.Collapse(0)\n$doiUrl = \"https://doi.org/10.1111/1365-2435.12799\"\n[void]$d.Hyperlinks.Add($linkRng, $doiUrl, \"\", \"\", $doiUrl)\n\n# Trailing period after the hyperlink.\n$p = $d.Paragraphs.Last\n$tail = $p.Range\n$tail.Collapse(0)                   # wdCollapseEnd\n$tail.InsertAfter(\".\")\n\n# Space between the page range and the hyperlink.\n$p = $d.Paragraphs.Last\n$r = $p.Range\n$r.Collapse(1)                      # wdCollapseStart\n$r.InsertBefore(\" \")\n\n# Volume/issue/page range, right before that space.\n$p = $d.Paragraphs.Last\n$r = $p.Range\n$r.Collapse(1)\n$r.InsertBefore(\"30 (12): 1878\u201382.\")\n\n# Space before the volume info.\n$p = $d.Paragraphs.Last\n$r = $p.Range\n$r.Collapse(1)\n$r.InsertBefore(\" \")\n\n# Italic journal title, right before that space.\n$p = $d.Paragraphs.Last\n$r = $p.Range\n$r.Collapse(1)\n$italicStart = $r.Start\n$journal = \"Functional Ecology\"\n$r.InsertBefore($journal)\n$italicRng = $d.Range($italicStart, $italicStart + $journal.Length)\n$italicRng.Font.Italic = 1\n\n# Space before the journal title.\n$p = $d.Paragraphs.Last\n$r = $p.Range\n$r.Collapse(1)\n$r.InsertBefore(\" \")\n\n# Leading citation text (authors, year, title) at the very start.\n$p = $d.Paragraphs.Last\n$r = $p.Range\n$r.Collapse(1)\n$lead = \"Poisot, Timoth\u00e9e, Daniel B. Stouffer, and Sonia K\u00e9fi. 2016. \u201cDescribe, Understand and Predict: Why Do We Need Networks in Ecology?\u201d\"\n$r.InsertBefore($lead)\n\n# --- 3. Wrap the new paragraph in bookmarks ---------------------------\n$p = $d.Paragraphs.Last\n$paraRange = $p.Range\n$d.Bookmarks.Add(\"ref-PoisStou16\", $paraRange)\n$d.Bookmarks.Add(\"refs\", $paraRange)\n"}
